$d = $word.ActiveDocument

# 1. Detection limit paragraph: remove CEBPA mention before ASXL1, move it later with TERT
$d.Content.Find.Execute(
    "with the exception of CEBPA (detection limit ~ 10%) and ASXL1 c.1934",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "with the exception of ASXL1 c.1934",
    2)

# 2. Update detection limit percentages and add CEBPA/TERT mention
$d.Content.Find.Execute(
    "(detection limit ~ 5%). This assay",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(detection limit ~ 5%-10%), CEBPA and TERT (detection limit ~ 10%). This assay",
    2)

# 3. Update VAF/CV% statistics
$d.Content.Find.Execute(
    "VAFs of 5%, 10%-20%, 30%-40% and 50% are on average, 10.2%, 10.4%, 3.5% and 4.4%, respectively",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "VAFs of 2%, 5%-10%, 20%-40% and 50% are on average, 15.4%, 8.6%, 4.0% and 1.8%, respectively",
    2)

# 4. Update saved date field text
$d.Content.Find.Execute(
    "30-Oct-2023",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "16-Nov-2023",
    2)
